$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.538.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -6.81%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.239.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -7.98%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '178.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -12.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '511.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.54%  '
$ws.Range("E7").Value = '  -1.36%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.236.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.612'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.98'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.97%  '
$ws.Range("E12").Value = '  -9.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000251'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.749.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.73%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.116'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.97%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.235.45'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '62.506.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -8.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.939'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '367.78'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.13'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -9.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.62%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.53%  '
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.75%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.18%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '28.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '628.09'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.90%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.27%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.10'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.75%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.105'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.35%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.90'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.88%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.09%  '
$ws.Range("B37").Value = 'TheGraph'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.390'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.06%  '
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.78'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -11.21%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.911.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.55%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.123'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.13%  '
$ws.Range("B42").Value = 'PEPE'
$ws.Range("C42").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₃0643'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.92%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.65%  '
$ws.Range("B44").Value = 'ThetaToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.64'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -14.30%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.58'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.36%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0385'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.90%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.69%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.124'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.02%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.26%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.40'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -15.49%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '128.95'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.00%  '
